$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that no longer hold data in the re-exported sheet
$ws.Range("D3").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("L5").ClearContents()

# Full re-exported data: column A now carries the running image id (19..140),
# rows 6-8, 10-12, 27, 39, 57, 63, 74, 80, 84, 91, 95, 97, 100, 106, 108, 114
# are newly populated (id only), and row 105 / new cols reflect updated counts.
$rowsData = @{
    2 = @{ A=19; L=1 }
    3 = @{ A=20; I=1; K=1; L=1 }
    4 = @{ A=21 }
    5 = @{ A=22; J=1 }
    6 = @{ A=23 }
    7 = @{ A=24 }
    8 = @{ A=25; D=1 }
    9 = @{ A=26; K=1 }
    10 = @{ A=27 }
    11 = @{ A=28 }
    12 = @{ A=29 }
    13 = @{ A=30; J=1 }
    14 = @{ A=31; I=1 }
    15 = @{ A=32; I=1 }
    16 = @{ A=33; J=1 }
    17 = @{ A=34; L=1 }
    18 = @{ A=35; J=1; K=1 }
    19 = @{ A=36; H=1; J=1 }
    20 = @{ A=37; E=1 }
    21 = @{ A=38; J=1 }
    22 = @{ A=39; K=1 }
    23 = @{ A=40; L=1 }
    24 = @{ A=41; J=1 }
    25 = @{ A=42; D=1 }
    26 = @{ A=43; K=1 }
    27 = @{ A=44 }
    28 = @{ A=45; L=1 }
    29 = @{ A=46; J=1 }
    30 = @{ A=47; L=1 }
    31 = @{ A=48; I=1; K=1 }
    32 = @{ A=49; I=1 }
    33 = @{ A=50; L=1 }
    34 = @{ A=51; K=1 }
    35 = @{ A=52; K=1 }
    36 = @{ A=53; L=1 }
    37 = @{ A=54; I=1; L=1 }
    38 = @{ A=55; L=1 }
    39 = @{ A=56 }
    40 = @{ A=57; F=1 }
    41 = @{ A=58; L=1 }
    42 = @{ A=59; K=1 }
    43 = @{ A=60; L=1 }
    44 = @{ A=61; L=1 }
    45 = @{ A=62; K=1 }
    46 = @{ A=63; J=1 }
    47 = @{ A=64; K=1 }
    48 = @{ A=65; K=1 }
    49 = @{ A=66; L=1 }
    50 = @{ A=67; K=1 }
    51 = @{ A=68; I=1 }
    52 = @{ A=69; J=1 }
    53 = @{ A=70; K=1 }
    54 = @{ A=71; L=1 }
    55 = @{ A=72; H=1 }
    56 = @{ A=73; I=1 }
    57 = @{ A=74 }
    58 = @{ A=75; G=1 }
    59 = @{ A=76; L=1 }
    60 = @{ A=77; E=1 }
    61 = @{ A=78; F=1 }
    62 = @{ A=79; K=1 }
    63 = @{ A=80 }
    64 = @{ A=81; L=1 }
    65 = @{ A=82; L=1 }
    66 = @{ A=83; L=1 }
    67 = @{ A=84; E=1 }
    68 = @{ A=85; J=1 }
    69 = @{ A=86; G=1 }
    70 = @{ A=87; H=1 }
    71 = @{ A=88; K=1 }
    72 = @{ A=89; I=1 }
    73 = @{ A=90; I=1 }
    74 = @{ A=91 }
    75 = @{ A=92; I=1 }
    76 = @{ A=93; I=1; K=1 }
    77 = @{ A=94; H=1 }
    78 = @{ A=95; J=1 }
    79 = @{ A=96; I=1 }
    80 = @{ A=97 }
    81 = @{ A=98; K=1 }
    82 = @{ A=99; L=1 }
    83 = @{ A=100; H=1 }
    84 = @{ A=101 }
    85 = @{ A=102; J=1 }
    86 = @{ A=103; G=1 }
    87 = @{ A=104; K=1 }
    88 = @{ A=105; H=1 }
    89 = @{ A=106; J=1; L=1 }
    90 = @{ A=107; G=1 }
    91 = @{ A=108 }
    92 = @{ A=109; I=1 }
    93 = @{ A=110; K=1 }
    94 = @{ A=111; K=1 }
    95 = @{ A=112 }
    96 = @{ A=113; I=1; L=2 }
    97 = @{ A=114 }
    98 = @{ A=115; H=1 }
    99 = @{ A=116; J=1 }
    100 = @{ A=117 }
    101 = @{ A=118; K=1 }
    102 = @{ A=119; K=1 }
    103 = @{ A=120; L=1 }
    104 = @{ A=121; L=1 }
    105 = @{ A=122; J=1; K=1 }
    106 = @{ A=123 }
    107 = @{ A=124; H=1 }
    108 = @{ A=125 }
    109 = @{ A=126; J=1 }
    110 = @{ A=127; L=1 }
    111 = @{ A=128; H=1 }
    112 = @{ A=129; J=1 }
    113 = @{ A=130; L=1 }
    114 = @{ A=131 }
    115 = @{ A=132; L=1 }
    116 = @{ A=133; C=1 }
    117 = @{ A=134; K=1 }
    118 = @{ A=135; K=1 }
    119 = @{ A=136; K=1 }
    120 = @{ A=137; K=1 }
    121 = @{ A=138; L=1 }
    122 = @{ A=139; K=1 }
    123 = @{ A=140 }
}

foreach ($r in $rowsData.Keys) {
    $rowVals = $rowsData[$r]
    foreach ($col in $rowVals.Keys) {
        $addr = "$col$r"
        $ws.Range($addr).Value = $rowVals[$col]
    }
}
